$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Majorelle Magdy, Dr. Eman Tantawi"
$ws.Range("G3").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Amira Sobhy"
$ws.Range("G4").Value = "Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Amira Sobhy"
$ws.Range("G5").Value = "Dr. Nesma, Dr. Nourhan Mahmoud, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Mohammad El-Tanany, Dr. Veronia Rafat, Dr. Hanan Ragab"
$ws.Range("G6").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Asmaa Reda, Dr. Nourhan Mahmoud, Dr. Eman Tantawi, Dr. Nahla Nagiub, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Amira Sobhy"
$ws.Range("G7").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Amira Sobhy"
$ws.Range("G8").Value = "Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Manar Montaser, Dr. Majorelle Magdy, Administrator, Dr. Shimaa Ahmad Mekki"
$ws.Range("G9").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Asmaa Reda, Dr. Hend Mahmoud, Dr. Majorelle Magdy, Dr. Manar Montaser, Dr. Rana Abo-Zaid, Dr. Gehan Adel, Dr. Amira Sobhy"
$ws.Range("G10").Value = "Dr. Alshimaa Atef, Dr. Heba Mahmoud Ali, Dr. Sara Wael, Dr. Servinaz Sayed Mohammad, Dr. Rana Abo-Zaid, Dr. Gehan Adel, Dr. Shimaa Ahmad Mekki"
$ws.Range("G11").Value = "Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Hend Mahmoud"
$ws.Range("G12").Value = "Dr. Salma El-Gendy, Administrator"
$ws.Range("G13").Value = "Dr. Shimaa Ashraf, Dr. Safa Hany, D Wessam Atef, Dr. Omnia Mohammad, Dr. Mariam Nour El-Din"
$ws.Range("G14").Value = "Dr. Shimaa Ashraf, Dr. Safa Hany"
$ws.Range("G17").Value = "Dr. Dina Adel, Dr. Basma Hamed, Dr. Nourhan Osama, Dr. Arwa Al-Sayed, Dr. Esraa Mostafa, Dr. Sarah Abdelmohsen, Dr. Madeha Saeed, Dr. Marwa Mustafa, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya"
$ws.Range("G19").Value = "Dr. Sarah Mahdy, D Mariam E. Mohammad"
$ws.Range("G22").Value = "Dr. Amr Saeed, Dr. Nancy Abd Al-Shafy"
$ws.Range("G24").Value = "Dr. Remon, Dr. Marina Atef, Dr. Youstina Magdy, Dr. Aya Emad, Dr. Yasmin, Dr. Maryam Ashraf, Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Neveen Nashaat"
$ws.Range("G25").Value = "Dr. Eman Samir Gabry, Dr. Remon, Dr. Marina Atef, Dr. Youstina Magdy, Dr. Aya Emad, Dr. Ola Abd Al-Fattah, Dr. Abdullah El-Agrody"
$ws.Range("G27").Value = "Dr. Remon, Dr. Eman Samir Gabry, Dr. Eman Mohammad Al, Dr. Yasmin, Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Ola Abd Al-Fattah, Dr. Neveen Nashaat"
$ws.Range("G28").Value = "Dr. Remon, Dr. Eman Samir Gabry, Dr. Nardine, Dr. Aya Hanafy, Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Neveen Nashaat, Dr. Abdullah El-Agrody"
$ws.Range("G29").Value = "Dr. Eman Samir Gabry, Dr. Remon, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah, Dr. Naema Gomaa, Dr. Monica"
$ws.Range("G30").Value = "Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Amira Sobhy"
$ws.Range("G31").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Amira Sobhy"
$ws.Range("G32").Value = "Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Amira Sobhy"
$ws.Range("G33").Value = "Dr. Nesma, Dr. Nourhan Mahmoud, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Mohammad El-Tanany, Dr. Veronia Rafat, Dr. Hanan Ragab"
$ws.Range("G34").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Asmaa Reda, Dr. Nourhan Mahmoud, Dr. Eman Tantawi, Dr. Nahla Nagiub, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Amira Sobhy"
$ws.Range("G35").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Amira Sobhy"
$ws.Range("G36").Value = "Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Manar Montaser, Dr. Majorelle Magdy, Administrator, Dr. Shimaa Ahmad Mekki"
$ws.Range("G37").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Asmaa Reda, Dr. Hend Mahmoud, Dr. Majorelle Magdy, Dr. Manar Montaser, Dr. Rana Abo-Zaid, Dr. Gehan Adel, Dr. Amira Sobhy"
$ws.Range("G38").Value = "Dr. Alshimaa Atef, Dr. Heba Mahmoud Ali, Dr. Sara Wael, Dr. Servinaz Sayed Mohammad, Dr. Rana Abo-Zaid, Dr. Gehan Adel, Dr. Shimaa Ahmad Mekki"
$ws.Range("G39").Value = "Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Hend Mahmoud"
$ws.Range("G40").Value = "Dr. Salma El-Gendy, Administrator"
$ws.Range("G41").Value = "Dr. Shimaa Ashraf, Dr. Safa Hany, D Wessam Atef, Dr. Omnia Mohammad, Dr. Mariam Nour El-Din"
$ws.Range("G42").Value = "Dr. Shimaa Ashraf, Dr. Safa Hany"
$ws.Range("G45").Value = "Dr. Dina Adel, Dr. Basma Hamed, Dr. Nourhan Osama, Dr. Arwa Al-Sayed, Dr. Esraa Mostafa, Dr. Sarah Abdelmohsen, Dr. Madeha Saeed, Dr. Marwa Mustafa, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya"
$ws.Range("G47").Value = "Dr. Sarah Mahdy, D Mariam E. Mohammad"
$ws.Range("G50").Value = "Dr. Amr Saeed, Dr. Nancy Abd Al-Shafy"
$ws.Range("G52").Value = "Dr. Remon, Dr. Marina Atef, Dr. Youstina Magdy, Dr. Aya Emad, Dr. Yasmin, Dr. Maryam Ashraf, Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Neveen Nashaat"
$ws.Range("G53").Value = "Dr. Eman Samir Gabry, Dr. Remon, Dr. Marina Atef, Dr. Youstina Magdy, Dr. Aya Emad, Dr. Ola Abd Al-Fattah, Dr. Abdullah El-Agrody"
$ws.Range("G55").Value = "Dr. Remon, Dr. Eman Samir Gabry, Dr. Eman Mohammad Al, Dr. Yasmin, Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Ola Abd Al-Fattah, Dr. Neveen Nashaat"
$ws.Range("G56").Value = "Dr. Remon, Dr. Eman Samir Gabry, Dr. Nardine, Dr. Aya Hanafy, Dr. Wafaa Ebida, Dr. Salma Hassan, Dr. Neveen Nashaat, Dr. Abdullah El-Agrody"
$ws.Range("G57").Value = "Dr. Eman Samir Gabry, Dr. Remon, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah, Dr. Naema Gomaa, Dr. Monica"
